# Insert a new weekly record at row 154 for
# "Hortaliza, Terminal La Palmera de La Serena - Cebollín", shifting the
# existing rows 154-215 down to 155-216.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(154).Insert()

$ws.Range("A154").Value = 8
$ws.Range("B154").Value = "Terminal La Palmera de La Serena"
$ws.Range("C154").Value = "Coquimbo"
$ws.Range("D154").Value = 44755
$ws.Range("E154").Value = 4
$ws.Range("F154").Value = 100112037
$ws.Range("G154").Value = "Cebollín"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 2000
$ws.Range("K154").Value = 1400
$ws.Range("L154").Value = 1600
$ws.Range("M154").Value = 1500
$ws.Range("N154").Value = "`$/paquete 6 unidades"
$ws.Range("O154").Value = "Provincia del Elquí"
$ws.Range("P154").Value = 250
$ws.Range("Q154").Value = 6
$ws.Range("R154").Value = "Hortaliza"
